# Fruta / hortaliza, semanal
# Reorders the data rows (2-18) of the sheet: each destination row receives
# the full record (columns A:R) that previously lived in a different source
# row, per the mapping below (rows 9 and 17 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 18
$firstCol = 1   # A
$lastCol = 18   # R

# destination row -> source row
$mapping = @{
    2  = 12
    3  = 13
    4  = 11
    5  = 10
    6  = 18
    7  = 3
    8  = 4
    9  = 9
    10 = 15
    11 = 16
    12 = 14
    13 = 2
    14 = 5
    15 = 8
    16 = 7
    17 = 17
    18 = 6
}

# Snapshot every source row's values before any writes happen, since this
# is a cyclic permutation (not simple swaps).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowValues[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowValues
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowValues = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowValues[$c]
    }
}
